$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 25-31 (this shifts old rows 32,33 up to 25,26)
$ws.Range("A25:A31").EntireRow.Delete()

# Delete rows 27-34 (removes the remaining data that had shifted up from rows 36-41)
$ws.Range("A27:A34").EntireRow.Delete()

# Update the selection to E20 (also clears the scrolled topLeftCell view state)
$ws.Range("E20").Select()
